# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) for a batch of leve rows across the crafting
# job sheets, as pulled from the latest market data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 47.9
$ws.Range("I11").Value = 47.9
$ws.Range("K11").Value = 47.9
$ws.Range("M11").Value = 92.09999999999999

$ws.Range("H31").Value = 1483.8334
$ws.Range("I31").Value = 380.6
$ws.Range("J31").Value = 7000
$ws.Range("K31").Value = 1141.8
$ws.Range("L31").Value = 21000
$ws.Range("M31").Value = -911.8000000000002
$ws.Range("N31").Value = -21460

$ws.Range("H40").Value = 1492.3125
$ws.Range("I40").Value = 1166.6666
$ws.Range("J40").Value = 1911
$ws.Range("K40").Value = 1166.6666
$ws.Range("L40").Value = 1911
$ws.Range("M40").Value = -991.6666
$ws.Range("N40").Value = -2261

$ws.Range("H64").Value = 3172.8948
$ws.Range("I64").Value = 2897.8333
$ws.Range("J64").Value = 3299.8462
$ws.Range("K64").Value = 2897.8333
$ws.Range("L64").Value = 3299.8462
$ws.Range("M64").Value = -2649.8333
$ws.Range("N64").Value = -3795.8462

$ws.Range("H67").Value = 3172.8948
$ws.Range("I67").Value = 2897.8333
$ws.Range("J67").Value = 3299.8462
$ws.Range("K67").Value = 2897.8333
$ws.Range("L67").Value = 3299.8462
$ws.Range("M67").Value = -2039.8333
$ws.Range("N67").Value = -5015.8462

$ws.Range("H87").Value = 42795
$ws.Range("J87").Value = 42795
$ws.Range("L87").Value = 42795
$ws.Range("N87").Value = -45291

$ws.Range("H90").Value = 42795
$ws.Range("J90").Value = 42795
$ws.Range("L90").Value = 128385
$ws.Range("N90").Value = -140865

$ws.Range("H129").Value = 1085.9038
$ws.Range("J129").Value = 1123.8163
$ws.Range("L129").Value = 3371.4489
$ws.Range("N129").Value = -13371.4489

$ws.Range("H132").Value = 2939.2727
$ws.Range("I132").Value = 2939.2727
$ws.Range("K132").Value = 8817.8181
$ws.Range("M132").Value = -6287.8181

$ws.Range("H133").Value = 50780
$ws.Range("J133").Value = 50780
$ws.Range("L133").Value = 50780
$ws.Range("N133").Value = -60900

$ws.Range("H137").Value = 1809.4667
$ws.Range("I137").Value = 1569.8572
$ws.Range("J137").Value = 2019.125
$ws.Range("K137").Value = 4709.571599999999
$ws.Range("L137").Value = 6057.375
$ws.Range("M137").Value = -2159.571599999999
$ws.Range("N137").Value = -11157.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1588.7646
$ws.Range("I61").Value = 1547.0714
$ws.Range("K61").Value = 1547.0714
$ws.Range("M61").Value = -1335.0714

$ws.Range("H97").Value = 76924080
$ws.Range("I97").Value = 824.8
$ws.Range("K97").Value = 824.8
$ws.Range("M97").Value = -328.8

$ws.Range("H122").Value = 2260.476
$ws.Range("I122").Value = 1526.1666
$ws.Range("J122").Value = 6666.3335
$ws.Range("K122").Value = 4578.4998
$ws.Range("L122").Value = 19999.0005
$ws.Range("M122").Value = -2128.4998
$ws.Range("N122").Value = -24899.0005

$ws.Range("H132").Value = 26530.715
$ws.Range("I132").Value = 2452.0625
$ws.Range("J132").Value = 103582.4
$ws.Range("K132").Value = 7356.1875
$ws.Range("L132").Value = 310747.2
$ws.Range("M132").Value = -4826.1875
$ws.Range("N132").Value = -315807.2

$ws.Range("H136").Value = 1588.7646
$ws.Range("I136").Value = 1547.0714
$ws.Range("K136").Value = 4641.2142
$ws.Range("M136").Value = -2091.2142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 575.9231
$ws.Range("I94").Value = 351.26315
$ws.Range("J94").Value = 1185.7142
$ws.Range("K94").Value = 351.26315
$ws.Range("L94").Value = 1185.7142
$ws.Range("M94").Value = 99.73685
$ws.Range("N94").Value = -2087.7142

$ws.Range("H107").Value = 950
$ws.Range("J107").Value = 1166.3334
$ws.Range("L107").Value = 1166.3334
$ws.Range("N107").Value = -5006.3334

$ws.Range("H134").Value = 3976.5518
$ws.Range("I134").Value = 4481.6665
$ws.Range("J134").Value = 1552
$ws.Range("K134").Value = 13444.9995
$ws.Range("L134").Value = 4656
$ws.Range("M134").Value = -10909.9995
$ws.Range("N134").Value = -9726

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1073.3334
$ws.Range("I16").Value = 810
$ws.Range("J16").Value = 1600
$ws.Range("K16").Value = 810
$ws.Range("L16").Value = 1600
$ws.Range("M16").Value = -523
$ws.Range("N16").Value = -2174

$ws.Range("H58").Value = 17627.645
$ws.Range("I58").Value = 1345.75
$ws.Range("J58").Value = 47231.09
$ws.Range("K58").Value = 1345.75
$ws.Range("L58").Value = 47231.09
$ws.Range("M58").Value = -1142.75
$ws.Range("N58").Value = -47637.09

$ws.Range("H105").Value = 886.5
$ws.Range("I105").Value = 602
$ws.Range("J105").Value = 1455.5
$ws.Range("K105").Value = 602
$ws.Range("L105").Value = 1455.5
$ws.Range("M105").Value = 1145
$ws.Range("N105").Value = -4949.5

$ws.Range("H107").Value = 1217.4
$ws.Range("I107").Value = 544.2727
$ws.Range("K107").Value = 544.2727
$ws.Range("M107").Value = 1375.7273

$ws.Range("H113").Value = 1073.3334
$ws.Range("I113").Value = 810
$ws.Range("J113").Value = 1600
$ws.Range("K113").Value = 810
$ws.Range("L113").Value = 1600
$ws.Range("M113").Value = 1360
$ws.Range("N113").Value = -5940

$ws.Range("H122").Value = 1569.2858
$ws.Range("I122").Value = 1767.4286
$ws.Range("K122").Value = 5302.2858
$ws.Range("M122").Value = -2852.2858

$ws.Range("H132").Value = 4228
$ws.Range("I132").Value = 3373.6
$ws.Range("J132").Value = 5448.5713
$ws.Range("K132").Value = 10120.8
$ws.Range("L132").Value = 16345.7139
$ws.Range("M132").Value = -7590.799999999999
$ws.Range("N132").Value = -21405.7139

$ws.Range("H136").Value = 17627.645
$ws.Range("I136").Value = 1345.75
$ws.Range("J136").Value = 47231.09
$ws.Range("K136").Value = 4037.25
$ws.Range("L136").Value = 141693.27
$ws.Range("M136").Value = -1487.25
$ws.Range("N136").Value = -146793.27

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 666.6667
$ws.Range("I49").Value = 500
$ws.Range("J49").Value = 1000
$ws.Range("K49").Value = 1500
$ws.Range("L49").Value = 3000
$ws.Range("M49").Value = -1344
$ws.Range("N49").Value = -3312

$ws.Range("H131").Value = 720.17
$ws.Range("J131").Value = 727.70105
$ws.Range("L131").Value = 2183.10315
$ws.Range("N131").Value = -12263.10315

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N9").ClearContents()
$ws.Range("H9").Value = 406
$ws.Range("I9").Value = 406
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 406
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -236

$ws.Range("H132").Value = 32197.295
$ws.Range("I132").Value = 3582.2
$ws.Range("K132").Value = 10746.6
$ws.Range("M132").Value = -8216.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5500
$ws.Range("I22").Value = 5500.5
$ws.Range("J22").Value = 5499.5
$ws.Range("K22").Value = 5500.5
$ws.Range("L22").Value = 5499.5
$ws.Range("M22").Value = -5205.5
$ws.Range("N22").Value = -6089.5

$ws.Range("H27").Value = 5500
$ws.Range("I27").Value = 5500.5
$ws.Range("J27").Value = 5499.5
$ws.Range("K27").Value = 5500.5
$ws.Range("L27").Value = 5499.5
$ws.Range("M27").Value = -5393.5
$ws.Range("N27").Value = -5713.5

$ws.Range("H100").Value = 2516.1667
$ws.Range("I100").Value = 2200
$ws.Range("J100").Value = 2579.4
$ws.Range("K100").Value = 2200
$ws.Range("L100").Value = 2579.4
$ws.Range("M100").Value = -1659
$ws.Range("N100").Value = -3661.4

$ws.Range("H122").Value = 756966
$ws.Range("I122").Value = 1785031.6
$ws.Range("K122").Value = 5355094.800000001
$ws.Range("M122").Value = -5352644.800000001

$ws.Range("H136").Value = 1089.2069
$ws.Range("I136").Value = 1053.9546
$ws.Range("J136").Value = 1200
$ws.Range("K136").Value = 3161.8638
$ws.Range("L136").Value = 3600
$ws.Range("M136").Value = -611.8638000000001
$ws.Range("N136").Value = -8700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 17500
$ws.Range("J54").Value = 17500
$ws.Range("L54").Value = 17500
$ws.Range("N54").Value = -18540

$ws.Range("H81").Value = 71430030
$ws.Range("I81").Value = 1342.8572
$ws.Range("J81").Value = 142858720
$ws.Range("K81").Value = 2685.7144
$ws.Range("L81").Value = 285717440
$ws.Range("M81").Value = -1624.7144
$ws.Range("N81").Value = -285719562

$ws.Range("H84").Value = 71430030
$ws.Range("I84").Value = 1342.8572
$ws.Range("J84").Value = 142858720
$ws.Range("K84").Value = 13428.572
$ws.Range("L84").Value = 1428587200
$ws.Range("M84").Value = -8124.572
$ws.Range("N84").Value = -1428597808

$ws.Range("H107").Value = 1101.5834
$ws.Range("I107").Value = 708.1667
$ws.Range("J107").Value = 1495
$ws.Range("K107").Value = 2124.5001
$ws.Range("L107").Value = 4485
$ws.Range("M107").Value = -204.5001000000002
$ws.Range("N107").Value = -8325

$ws.Range("H132").Value = 1825.3043
$ws.Range("I132").Value = 1245.091
$ws.Range("J132").Value = 2357.1667
$ws.Range("K132").Value = 3735.273
$ws.Range("L132").Value = 7071.500100000001
$ws.Range("M132").Value = -1205.273
$ws.Range("N132").Value = -12131.5001

$ws.Range("H136").Value = 26318530
$ws.Range("I136").Value = 34484070
$ws.Range("K136").Value = 103452210
$ws.Range("M136").Value = -103449660
